# Update cryptos list snapshot (prices / 1h volume change %) and fix the
# rank-11/12 ordering swap between Polygon and WrappedEther.
#
# Numeric-looking "Price" strings (e.g. "1.000", "0.06530") are written with
# a leading apostrophe so Excel keeps them as literal text instead of
# re-parsing them as numbers (which would silently drop trailing zeros /
# thousands-style dots, e.g. "1.000" -> 1). Percent strings in column E
# already contain surrounding spaces so Excel stores them as text as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.619.81"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.880.97"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'249.52"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.2933"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "'0.06530"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'21.93"
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").Value = "'0.07753"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "'96.83"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7390"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.880.74"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'5.244"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "'274.53"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "30.584.91"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'13.16"
$ws.Range("E18").Value = "  -3.62%  "
$ws.Range("D19").Value = "'0.000007534"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "2.128.95"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'5.324"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'6.234"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Value = "'9.211"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "'163.80"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "'18.85"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").Value = "'1.344"
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").Value = "'0.09700"
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("D31").Value = "'1.507"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'4.288"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "'4.143"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").Value = "'0.04861"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").Value = "'1.126"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'0.6981"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "'0.01903"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").Value = "'2.785"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").Value = "'6.306"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "'74.98"
$ws.Range("E41").Value = "  +7.02%  "
$ws.Range("D42").Value = "'2.021"
$ws.Range("E42").Value = "  +4.35%  "
$ws.Range("D43").Value = "'0.4246"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'102.35"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").Value = "'9.330"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "'7.050"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "'35.60"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").Value = "'914.24"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "'0.05744"
$ws.Range("E51").Value = "  +2.08%  "
